# Apply scheduled-runner price/profit updates across all Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 77.2
$ws.Range("I9").Value = 35.333332
$ws.Range("J9").Value = 140
$ws.Range("K9").Value = 35.333332
$ws.Range("L9").Value = 140
$ws.Range("M9").Value = 133.666668
$ws.Range("N9").Value = -478
$ws.Range("H19").Value = 1656.7858
$ws.Range("I19").Value = 1372.75
$ws.Range("J19").Value = 1770.4
$ws.Range("K19").Value = 1372.75
$ws.Range("L19").Value = 1770.4
$ws.Range("M19").Value = -1197.75
$ws.Range("N19").Value = -2120.4
$ws.Range("H43").Value = 971.75
$ws.Range("I43").Value = 1150
$ws.Range("J43").Value = 936.1
$ws.Range("K43").Value = 1150
$ws.Range("L43").Value = 936.1
$ws.Range("M43").Value = -1081
$ws.Range("N43").Value = -1074.1
$ws.Range("H51").Value = 3955.913
$ws.Range("J51").Value = 5998.5
$ws.Range("L51").Value = 5998.5
$ws.Range("N51").Value = -6966.5
$ws.Range("H62").Value = 5248.3
$ws.Range("I62").Value = 4668.3335
$ws.Range("J62").Value = 5496.857
$ws.Range("K62").Value = 4668.3335
$ws.Range("L62").Value = 5496.857
$ws.Range("M62").Value = -4044.3335
$ws.Range("N62").Value = -6744.857
$ws.Range("H65").Value = 5248.3
$ws.Range("I65").Value = 4668.3335
$ws.Range("J65").Value = 5496.857
$ws.Range("K65").Value = 23341.6675
$ws.Range("L65").Value = 27484.285
$ws.Range("M65").Value = -20221.6675
$ws.Range("N65").Value = -33724.285
$ws.Range("H80").Value = 4894.9
$ws.Range("J80").Value = 6414.722
$ws.Range("L80").Value = 19244.166
$ws.Range("N80").Value = -21240.166
$ws.Range("H83").Value = 4894.9
$ws.Range("J83").Value = 6414.722
$ws.Range("L83").Value = 57732.498
$ws.Range("N83").Value = -67716.49799999999
$ws.Range("H129").Value = 1033.2433
$ws.Range("I129").Value = 770
$ws.Range("J129").Value = 1094.6666
$ws.Range("K129").Value = 2310
$ws.Range("L129").Value = 3283.9998
$ws.Range("M129").Value = 2690
$ws.Range("N129").Value = -13283.9998
$ws.Range("H132").Value = 1658.4286
$ws.Range("I132").Value = 1834.8572
$ws.Range("J132").Value = 1129.1428
$ws.Range("K132").Value = 5504.571599999999
$ws.Range("L132").Value = 3387.4284
$ws.Range("M132").Value = -2974.571599999999
$ws.Range("N132").Value = -8447.428400000001
$ws.Range("H138").Value = 6252320
$ws.Range("I138").Value = 1637.6086
$ws.Range("J138").Value = 22226286
$ws.Range("K138").Value = 4912.825800000001
$ws.Range("L138").Value = 66678858
$ws.Range("M138").Value = 227.1741999999995
$ws.Range("N138").Value = -66689138

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1460.3334
$ws.Range("I45").Value = 1370.9231
$ws.Range("J45").Value = 1692.8
$ws.Range("K45").Value = 1370.9231
$ws.Range("L45").Value = 1692.8
$ws.Range("M45").Value = -993.9231
$ws.Range("N45").Value = -2446.8
$ws.Range("H61").Value = 1371.3846
$ws.Range("I61").Value = 1317.8334
$ws.Range("K61").Value = 1317.8334
$ws.Range("M61").Value = -1105.8334
$ws.Range("H88").Value = 3235.3333
$ws.Range("I88").Value = 2306
$ws.Range("J88").Value = 3700
$ws.Range("K88").Value = 2306
$ws.Range("L88").Value = 3700
$ws.Range("M88").Value = -1900
$ws.Range("N88").Value = -4512
$ws.Range("H91").Value = 3235.3333
$ws.Range("I91").Value = 2306
$ws.Range("J91").Value = 3700
$ws.Range("K91").Value = 2306
$ws.Range("L91").Value = 3700
$ws.Range("M91").Value = -902
$ws.Range("N91").Value = -6508
$ws.Range("H132").Value = 2684.8518
$ws.Range("I132").Value = 2039.6316
$ws.Range("J132").Value = 4217.25
$ws.Range("K132").Value = 6118.8948
$ws.Range("L132").Value = 12651.75
$ws.Range("M132").Value = -3588.8948
$ws.Range("N132").Value = -17711.75
$ws.Range("H136").Value = 1371.3846
$ws.Range("I136").Value = 1317.8334
$ws.Range("K136").Value = 3953.5002
$ws.Range("M136").Value = -1403.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1838.4348
$ws.Range("I86").Value = 1663.8823
$ws.Range("J86").Value = 2333
$ws.Range("K86").Value = 1663.8823
$ws.Range("L86").Value = 2333
$ws.Range("M86").Value = -540.8823
$ws.Range("N86").Value = -4579
$ws.Range("H89").Value = 1838.4348
$ws.Range("I89").Value = 1663.8823
$ws.Range("J89").Value = 2333
$ws.Range("K89").Value = 8319.4115
$ws.Range("L89").Value = 11665
$ws.Range("M89").Value = -2703.4115
$ws.Range("N89").Value = -22897
$ws.Range("H107").Value = 1387.55
$ws.Range("I107").Value = 1365.6875
$ws.Range("J107").Value = 1475
$ws.Range("K107").Value = 1365.6875
$ws.Range("L107").Value = 1475
$ws.Range("M107").Value = 554.3125
$ws.Range("N107").Value = -5315
$ws.Range("H137").Value = 36235.168
$ws.Range("J137").Value = 36235.168
$ws.Range("L137").Value = 36235.168
$ws.Range("N137").Value = -46435.168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1035.3784
$ws.Range("I31").Value = 969.375
$ws.Range("J31").Value = 1457.8
$ws.Range("K31").Value = 969.375
$ws.Range("L31").Value = 1457.8
$ws.Range("M31").Value = -674.375
$ws.Range("N31").Value = -2047.8
$ws.Range("H34").Value = 1035.3784
$ws.Range("I34").Value = 969.375
$ws.Range("J34").Value = 1457.8
$ws.Range("K34").Value = 969.375
$ws.Range("L34").Value = 1457.8
$ws.Range("M34").Value = -767.375
$ws.Range("N34").Value = -1861.8
$ws.Range("H99").Value = 4400
$ws.Range("I99").Value = 5000
$ws.Range("J99").Value = 4100
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 4100
$ws.Range("M99").Value = -3502
$ws.Range("N99").Value = -7096
$ws.Range("H126").Value = 4400
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 4100
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 12300
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -17240
$ws.Range("H134").Value = 2346
$ws.Range("I134").Value = 2363.32
$ws.Range("J134").Value = 2284.1428
$ws.Range("K134").Value = 7089.960000000001
$ws.Range("L134").Value = 6852.428400000001
$ws.Range("M134").Value = -4554.960000000001
$ws.Range("N134").Value = -11922.4284
$ws.Range("H140").Value = 39540
$ws.Range("J140").Value = 39540
$ws.Range("L140").Value = 39540
$ws.Range("N140").Value = -49900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4800
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4800
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 14400
$ws.Range("N76").Value = -15166
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 4800
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4800
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 14400
$ws.Range("N79").Value = -17052
$ws.Range("M79").ClearContents()
$ws.Range("H98").Value = 600
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 3841.6667
$ws.Range("J100").Value = 3841.6667
$ws.Range("L100").Value = 11525.0001
$ws.Range("N100").Value = -13147.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52.8125
$ws.Range("I2").Value = 20.75
$ws.Range("J2").Value = 149
$ws.Range("K2").Value = 20.75
$ws.Range("L2").Value = 149
$ws.Range("M2").Value = 92.25
$ws.Range("N2").Value = -375
$ws.Range("H102").Value = 3138
$ws.Range("I102").Value = 3500
$ws.Range("J102").Value = 2957
$ws.Range("K102").Value = 3500
$ws.Range("L102").Value = 2957
$ws.Range("M102").Value = -1878
$ws.Range("N102").Value = -6201
$ws.Range("H132").Value = 2768.2903
$ws.Range("I132").Value = 2382.647
$ws.Range("J132").Value = 3236.5715
$ws.Range("K132").Value = 7147.941
$ws.Range("L132").Value = 9709.7145
$ws.Range("M132").Value = -4617.941
$ws.Range("N132").Value = -14769.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3666.111
$ws.Range("I122").Value = 3082.8333
$ws.Range("K122").Value = 9248.499899999999
$ws.Range("M122").Value = -6798.499899999999
$ws.Range("H132").Value = 3507.258
$ws.Range("I132").Value = 3491.682
$ws.Range("K132").Value = 10475.046
$ws.Range("M132").Value = -7945.045999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 44242.75
$ws.Range("J41").Value = 10327
$ws.Range("L41").Value = 10327
$ws.Range("N41").Value = -11107
$ws.Range("H136").Value = 2877.8823
$ws.Range("I136").Value = 2920.111
$ws.Range("K136").Value = 8760.332999999999
$ws.Range("M136").Value = -6210.332999999999
